$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.671.68'
$ws.Range("E2").Value = '  +2.73%  '
$ws.Range("D3").Value = '2.971.08'
$ws.Range("E3").Value = '  +1.98%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''595.36'
$ws.Range("E5").Value = '  +0.94%  '
$ws.Range("D6").Value = '''144.40'
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D8").Value = '2.972.09'
$ws.Range("E8").Value = '  +2.03%  '
$ws.Range("D9").Value = '''0.503'
$ws.Range("E9").Value = '  -0.36%  '
$ws.Range("D10").Value = '''7.31'
$ws.Range("E10").Value = '  +6.02%  '
$ws.Range("D11").Value = '''0.144'
$ws.Range("E11").Value = '  +2.12%  '
$ws.Range("D12").Value = '''0.445'
$ws.Range("E12").Value = '  +1.39%  '
$ws.Range("D13").Value = '''0.0000235'
$ws.Range("E13").Value = '  +4.62%  '
$ws.Range("D14").Value = '''33.60'
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").Value = '3.460.52'
$ws.Range("E16").Value = '  +2.03%  '
$ws.Range("D17").Value = '62.482.46'
$ws.Range("E17").Value = '  +2.62%  '
$ws.Range("D18").Value = '''6.73'
$ws.Range("E18").Value = '  +0.64%  '
$ws.Range("D19").Value = '2.966.09'
$ws.Range("E19").Value = '  +2.00%  '
$ws.Range("D20").Value = '''442.05'
$ws.Range("E20").Value = '  +2.16%  '
$ws.Range("D21").Value = '''13.59'
$ws.Range("E21").Value = '  +1.67%  '
$ws.Range("D22").Value = '''0.676'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("E23").Value = '  +0.97%  '
$ws.Range("D24").Value = '''81.94'
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("D25").Value = '''10.85'
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = '''12.02'
$ws.Range("E26").Value = '  +1.87%  '
$ws.Range("B27").Value = 'Fetch.AI'
$ws.Range("C27").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D27").Value = '''2.15'
$ws.Range("E27").Value = '  -1.59%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").Value = '''2.61'
$ws.Range("E29").Value = '  +0.89%  '
$ws.Range("E30").Value = '  +0.61%  '
$ws.Range("E31").Value = '  -7.10%  '
$ws.Range("D32").Value = '''26.53'
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").Value = '''0.107'
$ws.Range("E33").Value = '  -1.00%  '
$ws.Range("D34").Value = '''1.00'
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("D35").Value = '0.0₃0881'
$ws.Range("E35").Value = '  +1.44%  '
$ws.Range("D36").Value = '''0.997'
$ws.Range("E36").Value = '  -0.82%  '
$ws.Range("D37").Value = '''5.64'
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("D38").Value = '''49.84'
$ws.Range("E38").Value = '  +0.31%  '
$ws.Range("E39").Value = '  -1.13%  '
$ws.Range("D40").Value = '''2.00'
$ws.Range("E40").Value = '  +1.42%  '
$ws.Range("D41").Value = '''8.64'
$ws.Range("E41").Value = '  +1.02%  '
$ws.Range("E42").Value = '  -2.29%  '
$ws.Range("D43").Value = '''0.281'
$ws.Range("E43").Value = '  -0.55%  '
$ws.Range("D44").Value = '''38.98'
$ws.Range("E44").Value = '  -4.73%  '
$ws.Range("D45").Value = '2.697.75'
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Value = '''369.02'
$ws.Range("E46").Value = '  -2.35%  '
$ws.Range("D47").Value = '''0.0341'
$ws.Range("E47").Value = '  -1.49%  '
$ws.Range("D48").Value = '''133.83'
$ws.Range("E48").Value = '  +0.04%  '
$ws.Range("D50").Value = '''23.27'
$ws.Range("E50").Value = '  -2.12%  '
$ws.Range("E51").Value = '  -0.84%  '
